$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 6).Value = 4.9
$ws.Cells.Item(2, 9).Value = 1.93
$ws.Cells.Item(2, 10).Value = 3.45
$ws.Cells.Item(2, 14).Value = 3.1
$ws.Cells.Item(2, 16).Value = 1.73
$ws.Cells.Item(2, 18).Value = 1.28
$ws.Cells.Item(2, 20).Value = 2.04
$ws.Cells.Item(2, 21).Value = 1.83
$ws.Cells.Item(2, 22).Value = 2.08
$ws.Cells.Item(2, 23).Value = 1.22
$ws.Cells.Item(2, 35).Value = 120
$ws.Cells.Item(3, 12).Value = 1.41
$ws.Cells.Item(3, 15).Value = 1.32
$ws.Cells.Item(3, 16).Value = 1.85
$ws.Cells.Item(3, 17).Value = 1.92
$ws.Cells.Item(3, 18).Value = 1.32
$ws.Cells.Item(3, 20).Value = 1.81
$ws.Cells.Item(3, 28).Value = 10
$ws.Cells.Item(3, 29).Value = 9.6
$ws.Cells.Item(3, 36).Value = 27
$ws.Cells.Item(3, 37).Value = 27
$ws.Cells.Item(3, 41).Value = 600
$ws.Cells.Item(4, 6).Value = 1.27
$ws.Cells.Item(4, 7).Value = 1.37
$ws.Cells.Item(4, 8).Value = 8.199999999999999
$ws.Cells.Item(4, 10).Value = 5.2
$ws.Cells.Item(4, 12).Value = 1.21
$ws.Cells.Item(4, 16).Value = 3.1
$ws.Cells.Item(4, 21).Value = 2.08
$ws.Cells.Item(4, 22).Value = 1.09
$ws.Cells.Item(4, 40).Value = 5.9
$ws.Cells.Item(5, 6).Value = 1.59
$ws.Cells.Item(5, 7).Value = 3.05
$ws.Cells.Item(5, 11).Value = 8
$ws.Cells.Item(5, 14).Value = 1.48
$ws.Cells.Item(5, 16).Value = 1.48
$ws.Cells.Item(5, 18).Value = 1.19
$ws.Cells.Item(5, 19).Value = 1.01
$ws.Cells.Item(5, 23).Value = 1.58
$ws.Cells.Item(5, 34).Value = 60
$ws.Cells.Item(6, 6).Value = 1.79
$ws.Cells.Item(6, 7).Value = 1.9
$ws.Cells.Item(6, 8).Value = 4.5
$ws.Cells.Item(6, 10).Value = 3.85
$ws.Cells.Item(6, 20).Value = 1.8
$ws.Cells.Item(6, 22).Value = 1.25
$ws.Cells.Item(6, 23).Value = 2.1
$ws.Cells.Item(6, 31).Value = 65
$ws.Cells.Item(6, 32).Value = 11.5
$ws.Cells.Item(6, 34).Value = 38
$ws.Cells.Item(6, 36).Value = 22
$ws.Cells.Item(7, 6).Value = 1.91
$ws.Cells.Item(7, 11).Value = 4.1
$ws.Cells.Item(7, 16).Value = 1.88
$ws.Cells.Item(7, 32).Value = 970
$ws.Cells.Item(7, 33).Value = 40
$ws.Cells.Item(8, 21).Value = 1.51
$ws.Cells.Item(8, 23).Value = 6
$ws.Cells.Item(9, 8).Value = 4.6
$ws.Cells.Item(9, 9).Value = 5
$ws.Cells.Item(9, 10).Value = 3.95
$ws.Cells.Item(9, 14).Value = 4.7
$ws.Cells.Item(9, 17).Value = 1.73
$ws.Cells.Item(9, 19).Value = 2.58
$ws.Cells.Item(9, 21).Value = 2.26
$ws.Cells.Item(9, 24).Value = 25
$ws.Cells.Item(9, 41).Value = 48
$ws.Cells.Item(10, 19).Value = 3.9
$ws.Cells.Item(11, 16).Value = 1.59
$ws.Cells.Item(11, 17).Value = 2.62
$ws.Cells.Item(11, 21).Value = 1.8
$ws.Cells.Item(11, 26).Value = 30
$ws.Cells.Item(11, 33).Value = 11
$ws.Cells.Item(12, 6).Value = 2.18
$ws.Cells.Item(12, 9).Value = 3.4
$ws.Cells.Item(12, 10).Value = 4
$ws.Cells.Item(12, 14).Value = 5.1
$ws.Cells.Item(12, 16).Value = 2.4
$ws.Cells.Item(12, 21).Value = 2.5
$ws.Cells.Item(12, 24).Value = 23
$ws.Cells.Item(12, 26).Value = 28
$ws.Cells.Item(12, 29).Value = 9.6
$ws.Cells.Item(12, 30).Value = 14.5
$ws.Cells.Item(12, 31).Value = 34
$ws.Cells.Item(12, 33).Value = 11
$ws.Cells.Item(12, 34).Value = 16
$ws.Cells.Item(12, 35).Value = 38
$ws.Cells.Item(12, 36).Value = 30
$ws.Cells.Item(12, 38).Value = 28
$ws.Cells.Item(12, 39).Value = 60
$ws.Cells.Item(12, 41).Value = 24
$ws.Cells.Item(13, 7).Value = 2.3
$ws.Cells.Item(13, 9).Value = 3.4
$ws.Cells.Item(13, 10).Value = 3.9
$ws.Cells.Item(13, 15).Value = 1.12
$ws.Cells.Item(13, 19).Value = 1.97
$ws.Cells.Item(13, 21).Value = 3.05
$ws.Cells.Item(13, 22).Value = 1.41
$ws.Cells.Item(13, 23).Value = 1.78
$ws.Cells.Item(13, 28).Value = 55
$ws.Cells.Item(13, 35).Value = 980
$ws.Cells.Item(13, 40).Value = 9.199999999999999
$ws.Cells.Item(14, 6).Value = 2.14
$ws.Cells.Item(14, 9).Value = 3.3
$ws.Cells.Item(14, 14).Value = 6.4
$ws.Cells.Item(14, 22).Value = 1.43
$ws.Cells.Item(14, 41).Value = 970
$ws.Cells.Item(15, 16).Value = 2.9
$ws.Cells.Item(15, 22).Value = 1.3
$ws.Cells.Item(17, 9).Value = 2.82
$ws.Cells.Item(17, 14).Value = 2.88
$ws.Cells.Item(17, 17).Value = 2.58
$ws.Cells.Item(17, 19).Value = 5
$ws.Cells.Item(17, 24).Value = 9.4
$ws.Cells.Item(17, 25).Value = 8.800000000000001
$ws.Cells.Item(17, 32).Value = 19
$ws.Cells.Item(17, 33).Value = 14.5
$ws.Cells.Item(17, 35).Value = 65
$ws.Cells.Item(17, 37).Value = 46
$ws.Cells.Item(17, 38).Value = 70
$ws.Cells.Item(17, 39).Value = 190
$ws.Cells.Item(18, 18).Value = 1.53
$ws.Cells.Item(18, 20).Value = 2.38
$ws.Cells.Item(18, 27).Value = 890
$ws.Cells.Item(18, 28).Value = 8.4
$ws.Cells.Item(18, 29).Value = 15
$ws.Cells.Item(18, 31).Value = 360
$ws.Cells.Item(18, 36).Value = 8.800000000000001
$ws.Cells.Item(18, 38).Value = 46
$ws.Cells.Item(18, 40).Value = 4.9
$ws.Cells.Item(18, 41).Value = 430
$ws.Cells.Item(19, 6).Value = 2.88
$ws.Cells.Item(19, 7).Value = 2.9
$ws.Cells.Item(19, 8).Value = 2.78
$ws.Cells.Item(19, 9).Value = 2.82
$ws.Cells.Item(19, 16).Value = 1.8
$ws.Cells.Item(19, 17).Value = 2.22
$ws.Cells.Item(19, 22).Value = 1.55
$ws.Cells.Item(19, 23).Value = 1.52
$ws.Cells.Item(19, 31).Value = 32
$ws.Cells.Item(19, 40).Value = 34
